$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4057
$ws.Range("J17").Value = 4057
$ws.Range("L17").Value = 12171
$ws.Range("N17").Value = -12507
$ws.Range("H33").Value = 242.12
$ws.Range("I33").Value = 251.09525
$ws.Range("K33").Value = 251.09525
$ws.Range("M33").Value = -22.09524999999999
$ws.Range("H99").Value = 372.76923
$ws.Range("J99").Value = 320
$ws.Range("L99").Value = 960
$ws.Range("N99").Value = -3956
$ws.Range("H101").Value = 2640.5
$ws.Range("I101").Value = 3004.8
$ws.Range("J101").Value = 2380.2856
$ws.Range("K101").Value = 9014.400000000001
$ws.Range("L101").Value = 7140.8568
$ws.Range("M101").Value = -7392.400000000001
$ws.Range("N101").Value = -10384.8568
$ws.Range("H125").Value = 3263.9092
$ws.Range("J125").Value = 4689.2
$ws.Range("L125").Value = 42202.8
$ws.Range("N125").Value = -47122.8
$ws.Range("H137").Value = 3756.6902
$ws.Range("I137").Value = 2341.3333
$ws.Range("J137").Value = 6710.478
$ws.Range("K137").Value = 7023.999899999999
$ws.Range("L137").Value = 20131.434
$ws.Range("M137").Value = -4473.999899999999
$ws.Range("N137").Value = -25231.434
$ws.Range("H138").Value = 4662.9507
$ws.Range("I138").Value = 4431.1665
$ws.Range("J138").Value = 4688.2363
$ws.Range("K138").Value = 13293.4995
$ws.Range("L138").Value = 14064.7089
$ws.Range("M138").Value = -8153.499500000002
$ws.Range("N138").Value = -24344.7089
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5007.02
$ws.Range("I32").Value = 2609.0588
$ws.Range("J32").Value = 18595.467
$ws.Range("K32").Value = 2609.0588
$ws.Range("L32").Value = 18595.467
$ws.Range("M32").Value = -2322.0588
$ws.Range("N32").Value = -19169.467
$ws.Range("H60").Value = 16338
$ws.Range("I60").Value = 16338
$ws.Range("K60").Value = 16338
$ws.Range("M60").Value = -15605
$ws.Range("H61").Value = 10802.5
$ws.Range("I61").Value = 10175
$ws.Range("J61").Value = 11011.667
$ws.Range("K61").Value = 10175
$ws.Range("L61").Value = 11011.667
$ws.Range("M61").Value = -9963
$ws.Range("N61").Value = -11435.667
$ws.Range("H74").Value = 389352.8
$ws.Range("I74").Value = 716441.6
$ws.Range("K74").Value = 716441.6
$ws.Range("M74").Value = -715567.6
$ws.Range("H77").Value = 389352.8
$ws.Range("I77").Value = 716441.6
$ws.Range("K77").Value = 3582208
$ws.Range("M77").Value = -3577840
$ws.Range("H105").Value = 97071.67
$ws.Range("J105").Value = 97071.67
$ws.Range("L105").Value = 97071.67
$ws.Range("N105").Value = -104059.67
$ws.Range("H132").Value = 7582.7354
$ws.Range("I132").Value = 5799.9614
$ws.Range("K132").Value = 17399.8842
$ws.Range("M132").Value = -14869.8842
$ws.Range("H136").Value = 10802.5
$ws.Range("I136").Value = 10175
$ws.Range("J136").Value = 11011.667
$ws.Range("K136").Value = 30525
$ws.Range("L136").Value = 33035.001
$ws.Range("M136").Value = -27975
$ws.Range("N136").Value = -38135.001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 100002800
$ws.Range("I105").Value = 142859140
$ws.Range("K105").Value = 142859140
$ws.Range("M105").Value = -142857393
$ws.Range("H132").Value = 99804.5
$ws.Range("J132").Value = 99804.5
$ws.Range("L132").Value = 99804.5
$ws.Range("N132").Value = -109924.5
$ws.Range("H134").Value = 4244.7407
$ws.Range("I134").Value = 3004.9092
$ws.Range("K134").Value = 9014.7276
$ws.Range("M134").Value = -6479.7276
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28575904
$ws.Range("I31").Value = 45456780
$ws.Range("J31").Value = 8266.154
$ws.Range("K31").Value = 45456780
$ws.Range("L31").Value = 8266.154
$ws.Range("M31").Value = -45456485
$ws.Range("N31").Value = -8856.154
$ws.Range("H34").Value = 28575904
$ws.Range("I34").Value = 45456780
$ws.Range("J34").Value = 8266.154
$ws.Range("K34").Value = 45456780
$ws.Range("L34").Value = 8266.154
$ws.Range("M34").Value = -45456578
$ws.Range("N34").Value = -8670.154
$ws.Range("H132").Value = 18373.072
$ws.Range("I132").Value = 4567.4
$ws.Range("J132").Value = 39944.438
$ws.Range("K132").Value = 13702.2
$ws.Range("L132").Value = 119833.314
$ws.Range("M132").Value = -11172.2
$ws.Range("N132").Value = -124893.314
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 990.5
$ws.Range("J5").Value = 1951.25
$ws.Range("L5").Value = 5853.75
$ws.Range("N5").Value = -6077.75
$ws.Range("H34").Value = 1459
$ws.Range("I34").Value = 541.5454999999999
$ws.Range("J34").Value = 2089.75
$ws.Range("K34").Value = 1624.6365
$ws.Range("L34").Value = 6269.25
$ws.Range("M34").Value = -1540.6365
$ws.Range("N34").Value = -6437.25
$ws.Range("H95").Value = 8000
$ws.Range("J95").Value = 8000
$ws.Range("L95").Value = 24000
$ws.Range("N95").Value = -28118
$ws.Range("H135").Value = 990.5
$ws.Range("J135").Value = 1951.25
$ws.Range("L135").Value = 17561.25
$ws.Range("N135").Value = -22631.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 37735.855
$ws.Range("I46").Value = 9589.25
$ws.Range("J46").Value = 75264.664
$ws.Range("K46").Value = 9589.25
$ws.Range("L46").Value = 75264.664
$ws.Range("M46").Value = -9433.25
$ws.Range("N46").Value = -75576.664
$ws.Range("H122").Value = 5768.815
$ws.Range("I122").Value = 4963.0713
$ws.Range("J122").Value = 6636.5386
$ws.Range("K122").Value = 14889.2139
$ws.Range("L122").Value = 19909.6158
$ws.Range("M122").Value = -12439.2139
$ws.Range("N122").Value = -24809.6158
$ws.Range("H132").Value = 3190.8867
$ws.Range("I132").Value = 2619.0244
$ws.Range("K132").Value = 7857.073199999999
$ws.Range("M132").Value = -5327.073199999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4065.5483
$ws.Range("I22").Value = 2670
$ws.Range("J22").Value = 5073.4443
$ws.Range("K22").Value = 2670
$ws.Range("L22").Value = 5073.4443
$ws.Range("M22").Value = -2375
$ws.Range("N22").Value = -5663.4443
$ws.Range("H27").Value = 4065.5483
$ws.Range("I27").Value = 2670
$ws.Range("J27").Value = 5073.4443
$ws.Range("K27").Value = 2670
$ws.Range("L27").Value = 5073.4443
$ws.Range("M27").Value = -2563
$ws.Range("N27").Value = -5287.4443
$ws.Range("H46").Value = 7215.759
$ws.Range("I46").Value = 2477.5557
$ws.Range("K46").Value = 2477.5557
$ws.Range("M46").Value = -2289.5557
$ws.Range("H68").Value = 4172.4165
$ws.Range("I68").Value = 3497
$ws.Range("K68").Value = 3497
$ws.Range("M68").Value = -2748
$ws.Range("H71").Value = 4172.4165
$ws.Range("I71").Value = 3497
$ws.Range("K71").Value = 17485
$ws.Range("M71").Value = -13741
$ws.Range("H132").Value = 5576.8296
$ws.Range("I132").Value = 4964.185
$ws.Range("K132").Value = 14892.555
$ws.Range("M132").Value = -12362.555
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 11977.5
$ws.Range("J11").Value = 11977.5
$ws.Range("L11").Value = 11977.5
$ws.Range("N11").Value = -12261.5
$ws.Range("H107").Value = 1134.8077
$ws.Range("I107").Value = 1057.1818
$ws.Range("K107").Value = 3171.5454
$ws.Range("M107").Value = -1251.5454
$ws.Range("H122").Value = 3134.3462
$ws.Range("I122").Value = 2817.8635
$ws.Range("K122").Value = 8453.5905
$ws.Range("M122").Value = -6003.5905
$ws.Range("H132").Value = 8034.769
$ws.Range("I132").Value = 5628.857
$ws.Range("J132").Value = 10841.667
$ws.Range("K132").Value = 16886.571
$ws.Range("L132").Value = 32525.001
$ws.Range("M132").Value = -14356.571
$ws.Range("N132").Value = -37585.001
$ws.Range("H136").Value = 3403573
$ws.Range("I136").Value = 4609969.5
$ws.Range("J136").Value = 3728.5454
$ws.Range("K136").Value = 13829908.5
$ws.Range("L136").Value = 11185.6362
$ws.Range("M136").Value = -13827358.5
$ws.Range("N136").Value = -16285.6362
